# Config.xlsx update: add ArriveNow_CredentialName to Constants sheet,
# refresh the Assets sheet (drop the ArriveNowCredentials/EmailAddress
# rows and append the new reporting-related assets).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Constants sheet: new row 6 - ArriveNow_CredentialName
# ---------------------------------------------------------------
$wsConstants = $wb.Worksheets.Item("Constants")
$wsConstants.Range("A6").Value = "ArriveNow_CredentialName"
$wsConstants.Range("B6").Value = "12_ArriveNowCredentials"
$wsConstants.Range("C6").Value = "The name of the orchestrator asset where the ArriveNow credentials are stored"

# ---------------------------------------------------------------
# Assets sheet: rewrite rows 2-16 with the refreshed asset list
# ---------------------------------------------------------------
$wsAssets = $wb.Worksheets.Item("Assets")

$assetRows = @(
    @("ArriveNowURL", "12_ArriveNowURL", "Shared", "URL for ArriveNow Portal"),
    @("ArrivePortalURL", "12_ArrivePortalURL", "Shared", "URL for Arrive Portal"),
    @("ArriveTruckEntryURL", "12_ArriveTruckEntryURL", "Shared", "URL for ArriveNow Truck Entry Portal"),
    @("GDriveReportFolder", "12_GDriveReportFolder", "Shared", "ID for G Drive folder where reports are stored"),
    @("ToEmail", "12_ToEmail", "Shared", "Email addresses where the emails are going to be sent to"),
    @("CCEmail", "12_CCEmail", "Shared", "Email addresses copied to the emails that are going to be sent"),
    @("ArriveCarrierSearchURL", "12_ArriveCarrierSearchURL", "Shared", "URL for ArriveNow Carrier Search Portal"),
    @("ReportFileID", "12_ReportFileID", "Shared", "ID for G Sheet used to report execution outputs"),
    @("ReportCreationDate", "12_ReportCreationDate", "Shared", "Date used to check when the last report was created"),
    @("ReportFileURL", "12_ReportFileURL", "Shared", "URL for G Sheet used to report execution outputs"),
    @("CompletedCasesCount", "12_CompletedCasesCount", "Shared", "This asset holds the number of completed cases for the day"),
    @("ExceptionCasesCount", "12_ExceptionCasesCount", "Shared", "This asset holds the number of exception cases for the day"),
    @("MasterReportID", "12_MasterReportID", "Shared", "ID for G Sheet for the Master execution report"),
    @("MasterReportURL", "12_MasterReportURL", "Shared", "URL for G Sheet for the Master execution report"),
    @("TotalCasesCount", "12_TotalCasesCount", "Shared", "This asset holds the number of cases processed for the day")
)

$row = 2
foreach ($entry in $assetRows) {
    $wsAssets.Cells.Item($row, 1).Value = $entry[0]
    $wsAssets.Cells.Item($row, 2).Value = $entry[1]
    $wsAssets.Cells.Item($row, 3).Value = $entry[2]
    $wsAssets.Cells.Item($row, 4).Value = $entry[3]
    $row = $row + 1
}

# Remove the two now-unused placeholder rows from the bottom of Assets
$wsAssets.Rows.Item(999).EntireRow.Delete()
$wsAssets.Rows.Item(999).EntireRow.Delete()

# ---------------------------------------------------------------
# Selection / active-sheet state
# ---------------------------------------------------------------
$wsAssets.Range("A2:XFD2").Select()
$wsConstants.Activate()
$wsConstants.Range("C9").Select()
